$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "id" header at A1 (existing B1..G1 headers remain as-is)
$ws.Range("A1").Value = "id"

# New header columns H1, I1
$ws.Range("H1").Value = "uniprot_rxn_1"
$ws.Range("I1").Value = "uniprot_rxn_2"

# Copy header style (bold, centered, bordered) from B1 onto the new header cells
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1:I1").PasteSpecial(-4122) | Out-Null

# Row 2 updates
$ws.Range("A2").Value = 7
$ws.Range("D2").Value = 55.00491922689314
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "Q7WKM5, A3Q5P5, A1SMQ5, Q2S3H9, A1BE51, Q5LCW4, Q2YXM1, P53593, Q74EA5, Q5L593, B7LAD4, A7N158, P53587, B8GYZ1, B1YIL7, Q66DA0, P0A836, C1EP60, B6IZ97, P80865, B0BXJ5, A4WL16, A6T6F6, Q04QS3, O05966, A1VYP2, O82662, Q2N695, O28098, B1ZP17, Q7W789, B4TBD7, B0SAP0, Q9AB94, Q2JFA7, P0AGE9, Q0BKS4, Q0AKU3, C3L791, Q3SVK4, Q92I22, B3TZD9, B8GCC9, Q87A98, A8FD73, B1MJJ0, Q9CNZ0, B2VBR6, Q57RL3, A3NRX7, A8FKV6, Q4ZUW7, A4T6V7, A4TNT8, B0RPQ9, B5EZG1, Q6D7G2, B1JG57, A1UL96, Q5NMF5, B4RCH3, Q8F746, A0LIY8, A3DLR9, B1W3P6, A9L532, B2ST87, Q57AX2, Q5WZ04, P53312, B3EFL3, B2AHC1, Q4FP29, Q1AWI8, Q3ATR4, B4RL77, Q819X0, B9IVC5, Q6F8L4, Q2RV33, B4S612, A1WD26, O28733, B3CTW0, Q5E6M6, B7NMS9, P66871, Q1D6I9, Q883Z4, Q7MMN3, C3P5N7, B0BTV5, A8ZRW7, B4RYG2, A1RTB8, Q9AKE1, Q2T0U7, A6QGE5, Q21RF9, P0A839, B9LB77, Q6HEX9, A5U0Z1, B3QPP2, W0PFR9, A9WKU9, B2SDM1, B3EN95, P68209, B4ESR1, Q0A5K0, Q2GJ42, B4ST72, Q5PCM7, Q63WW2, Q9KA20, Q582T6, A7NNW3, Q8EFN8, Q3A7Y4, B3PQ90, Q8ZH00, Q6FYD1, P9WGC6, Q3IZ84, A8GUV0, A9IZV3, Q2GFL7, Q9JUT0, B2V8I4, A3MH60, Q255N7, Q12M63, Q32IK3, P0A838, A5VSQ3, Q3JVL5, B8ZU00, B1LLG1, P53400, P45101, Q02K73, Q3BQ83, B5FC28, A4SND2, B0T3D1, Q5X7K6, B8HDC2, Q2IIC1, A5CC52, Q65JP0, A9MTQ2, P9WGC7, B2HXG0, A9R2F2, B0VEF2, C1DM53, A0KJK9, A0RHK9, Q13DR5, A4W879, C3PNF9, Q324I4, Q8Y1Y3, A1SWX6, A0R3M4, Q9Z5H8, A7ZJA8, Q8NX01, C1ALT7, Q5NHF3, C0QLM0, A9IKE3, Q13UH9, A3N693, P66866, C3LP35, Q8A9M7, A4XBE5, Q5YPP8, C4K7D7, A5WC33, Q5F878, A3D5J3, A9N8R8, Q9P2R7, Q87RE8, A9VT74, Q3KFU6, Q6A6Y0, A3MWQ7, B9MHW0, A6GY30, Q6MBM7, Q9JZP4, A6WPA4, B9JTS8, Q2GEF1, Q6LTB5, O87840, B7HLG9, O28097, A1R3M0, Q2KYF7, Q1RH56, A6SUY3, A7X1L8, A9M4F8, A1S5I0, B0CIT0, P53399, B9L5I6, Q92I21, Q0I3A8, C0ZFB3, P99507, A6UDP2, A3QDH6, Q3B5J6, Q5QU21, A4IM83, A5IH21, A8GS29, Q475M1, B7GXK7, Q5WFP4, A7ZXY8, Q6GHI9, B5EMG2, Q8UC60, Q9HPP1, B1Y9F9, Q3SKM1, P53598, B6I7Z9, C4L8V5, Q0T6W6, Q1MAW1, Q65SV1, B2FRS3, A1VUC0, B6EHV6, Q1LTZ4, A6TKJ3, A6U164, P99071, B0TZ23, A7HT39, Q9P567, A6VZ30, Q2FHJ3, A0JT43, Q4L5U8, B3QWF8, C0PWE7, B7I6T2, Q4ULQ7, P53588, Q7UKI3, Q8ZVF3, B5BC72, B9KNB3, Q482S1, Q81WK0, A9KBQ4, B2UMH5, P45102, A1JRB6, B5QWG8, A8EYY8, Q21IW6, Q5L0N7, C0RFH1, A8A8L8, B0BAS9, B1XRZ4, O08371, P9WGC5, Q5HGI7, O27115, C5BL83, P80886, Q9EYG9, Q4UXU0, A4IZC2, C1CV60, Q5GSB6, B4UDP4, W0PAN5, C0R4N4, Q03184, B6JE33, O67547, A8M469, Q7M9A6, Q3YSV6, Q39D96, A9MJL9, Q8P676, B4EBY6, Q8CPH5, Q8LAD2, Q6ZL94, A0RPE8, A1UQV9, A0QBP0, Q0KE75, P66872, A4SZY0, Q7VVU5, B5ZSR9, A1RIT3, B1IY02, Q053V6, A0LZD4, B1X6Q8, A6VPM2, Q9PAH1, Q6AD57, B9LRR0, Q8NMK7, Q7N6V5, B4SZE2, A6Q7S7, A7MQX5, B3E7B3, Q0HW00, Q1B3V2, Q8DFQ1, Q9Z2I9, P53592, Q5HC25, C5BEM7, Q4ULQ8, B0B950, Q2G3G8, C4K0K1, Q73FW3, O67546, Q2W063, Q2SD35, B0KNW8, A6WXE8, A7FKR4, C0QQH3, O26663, Q1CFM0, Q2A253, A1V6W8, P0AGF1, A1AVI3, A4JHJ8, A6V7K5, A1WZ92, C5C916, Q4FVH9, Q7Z940, A0PWD1, Q9KQB5, B5XZD1, B8FLW6, B1YVY6, B3GX29, Q9P727, B8E770, A5EX54, Q0VPF8, C5BZT2, Q1I7L3, C5D8V2, Q28U74, B2I9R4, A8MM16, A4G230, B9JCF4, B1JAV3, B2K8F1, Q6GHJ0, Q21CW8, A0LRU4, A4VKQ1, Q0C5F3, C1A3F2, Q6MGW1, B2SYE8, P99070, O28732, A3M887, A5VB78, Q98EC5, A7NDR1, Q57663, Q732N4, B0VSL0, A4Y7Q5, Q122H7, C4ZWK2, Q1IZE7, Q7Z941, Q1IJS5, Q0VPF7, A9GD65, A7GYI5, B8F4Q2, Q148D5, Q3IGW5, Q1QQR3, B8D5I6, A8GB83, B7IUJ5, Q1RH55, A8F1K7, A0KVX6, A5F6G6, Q6ARL2, A4FP60, B8GN38, C6DCD6, A3PN13, A5W114, P0AGF0, Q88FB2, P66868, Q9RUY3, Q3Z476, Q15UW6, B7M5P1, Q07UX6, C1AWW6, P53591, B7J4G8, C3K6N0, A1KTM6, B1KPH0, Q55AI5, Q5HGI6, B0UUF3, Q11CV8, A1TEP3, A9AFT3, B2GM49, B3FHT4, B2HED1, B0TJP8, Q14IV5, P0A837, Q5HVN3, Q6K9N6, C1A8X8, Q5PBI8, B6IP72, Q5FH54, A7GRG6, Q47AN8, B7HDV9, A9B0I2, A1APQ8, Q30RQ3, B0R5U5, P53401, B7VI35, Q2FZ37, B4SEH6, Q021N6, Q2J3G9, A3MZH0, Q7U0Z0, Q7U0Z1, A8Z3S5, Q8PHL5, Q1H031, O13750, Q49X32, A5V080, Q2P5E6, B9DPG0, B8CR00, P66867, A8FX46, C1DB55, Q5HPU5, B3CME9, B6J8N7, A8H3H6, Q2K3F0, Q31FV4, B2TUB2, A1KH85, Q6NDB7, Q8GTQ9, C4L610, Q4QLA3, Q5P7A7, P86226, Q1BTW6, Q6G1M1, A7H4K9, A9M8R2, B1JY14, A0KA72, Q636J0, Q2NUM2, Q7NZ47, A6Q391, A7I2F2, A0Q5H4, Q5H2H4, Q68WU5, A5UCR4, A4WNM6, A5FHX8, Q5HPU4, A1WLG5, Q8CPH4, Q84LB6, B0U5F8, B5R689, B2S880, B3FHP0, Q64U25, A1TZS9, B8JAE0, A4XV90, A7Z4M9, Q3KKL8, Q9Z6T6, A1TW81, A9BQA5, Q0S533, A4YKC6, A8GNF5, Q822A1, Q9PL99, B5YQR7, Q4R517, B4U901, O94415, Q1CAG1, B3TZD8, B2UGE2, Q1QEK5, Q6G9W7, Q9YI37, A1AZH7, Q6DQL1, Q58643, A2S8U3, Q4KFY6, B7UVD3, Q48K68, Q89X60, Q72PA2, Q1KSE5, A4SDL6, A5CU44, Q080X8, Q62MG7, P9WGC4, P66870, A5UIW0, P66869, Q9KY56, Q0BC71, Q5V4Q4, B5FNF9, Q6G9W8, A5ISD0, B3Q760, B7N9W6, A9WJ74, Q0HJQ4, Q8KFE7, O97580, A5E936, A8LAW3, Q51567, C1F2F2, A1KAU3, Q3INN3, B4TQ53, P53597, A5FYZ9, O84828, Q9PHY1, B9EBC7"
$ws.Range("I2").Value = "Q8WVX9"

# Row 3 updates
$ws.Range("A3").Value = 15
$ws.Range("D3").Value = 9.351730564334575
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "Q9JID6, P33121, P97524, O14975, P18163, P41216, O35488"
$ws.Range("I3").Value = "Q8WVX9"

# Row 4 - new row
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "2mg"
$ws.Range("C4").Value = "mvacid"
$ws.Range("D4").Value = 45.78995244029403
$ws.Range("E4").Value = 0.3553921568627451
$ws.Range("F4").Value = 0.7107843137254901
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = "B3EY95"
$ws.Range("I4").Value = "Q8WVX9"

# Copy the "id" column header/number style (style index 1) from A3 down to A4
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null

$wb.Save()
